$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 60
$ws1.Range("F3").Value = 21599
$ws1.Range("F4").Value = 820
$ws1.Range("F6").Value = 1136
$ws1.Range("F8").Value = 8032
$ws1.Range("F9").Value = 561
$ws1.Range("F12").Value = 323
$ws1.Range("F13").Value = 71
$ws1.Range("F14").Value = 193
$ws1.Range("F15").Value = 185
$ws1.Range("F17").Value = 236
$ws1.Range("F19").Value = 1366
$ws1.Range("F20").Value = 562
$ws1.Range("F23").Value = 56
$ws1.Range("F24").Value = 93
$ws1.Range("F26").Value = 359
$ws1.Range("F27").Value = 1206
$ws1.Range("F28").Value = 65
$ws1.Range("F30").Value = 233
$ws1.Range("F33").Value = 155
$ws1.Range("F34").Value = 5143
$ws1.Range("F35").Value = 37
$ws1.Range("F37").Value = 59
$ws1.Range("F39").Value = 13277
$ws1.Range("F40").Value = 1374
$ws1.Range("F44").Value = 330
$ws1.Range("F45").Value = 456
$ws1.Range("F46").Value = 4075
$ws1.Range("F47").Value = 29

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 60
$ws4.Range("F3").Value = 21599
$ws4.Range("F4").Value = 1136
$ws4.Range("F6").Value = 8032
$ws4.Range("F7").Value = 561
$ws4.Range("F10").Value = 323
$ws4.Range("F11").Value = 71
$ws4.Range("F12").Value = 193
$ws4.Range("F13").Value = 185
$ws4.Range("F15").Value = 236
$ws4.Range("F17").Value = 1366
$ws4.Range("F18").Value = 562
$ws4.Range("F21").Value = 56
$ws4.Range("F22").Value = 93
$ws4.Range("F24").Value = 359
$ws4.Range("F25").Value = 1206
$ws4.Range("F26").Value = 65
$ws4.Range("F28").Value = 233
$ws4.Range("F33").Value = 155
$ws4.Range("F35").Value = 5143
$ws4.Range("F36").Value = 37
$ws4.Range("F38").Value = 59
$ws4.Range("F40").Value = 13277
$ws4.Range("F41").Value = 1374
$ws4.Range("F44").Value = 330
$ws4.Range("F45").Value = 456
$ws4.Range("F46").Value = 4075
$ws4.Range("F47").Value = 29
